$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.561537265777588
$ws.Range("B1").Value = 1.77750551700592
$ws.Range("C1").Value = 4.13076639175415
$ws.Range("D1").Value = 1.921987533569336
$ws.Range("E1").Value = 0.7990439534187317
